$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the 2023-06 row (row 38) of expense/income data that was previously blank.
$ws.Range("C38").Value = 270.39999999999998
$ws.Range("D38").Value = 203
$ws.Range("E38").Value = 11
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 3294.04
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Formula = "=(H38+I38)-(C38+D38+E38+F38+G38)"

# Move the active selection to where the user ended up after editing.
$ws.Range("G39").Select()
